$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 16:22"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 534494
$ws.Range("C4").Value = 1615
$ws.Range("D4").Value = 30548
$ws.Range("E4").Value = 483309
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 20637

# Row 17: Brasil
$ws.Range("B17").Value = 21040
$ws.Range("C17").Value = 78
$ws.Range("E17").Value = 19723
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 1144

# Row 74: Bosnia y Herzegovina
$ws.Range("B74").Value = 1000
$ws.Range("C74").Value = 54
$ws.Range("E74").Value = 769

# Row 82: Bulgaria
$ws.Range("B82").Value = 675
$ws.Range("C82").Value = 14
$ws.Range("E82").Value = 578

# Row 153: Zambia
$ws.Range("B153").Value = 43
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 30
$ws.Range("E153").Value = 11
